# Updates odds values in Sheet1 to reflect the latest FlashScore scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "O3" = 1.91
    "P3" = 1.8
    "Q3" = 3.05
    "R3" = 1.38
    "S3" = 4.2
    "T3" = 1.22
    "U3" = 8
    "V3" = 1.08
    "Y3" = 1.93
    "Z3" = 1.88
    "AA13" = 1.65
    "AB13" = 1.98
    "AC13" = 14.5
    "AI13" = 12.5
    "AJ13" = 7.6
    "AM13" = 7.9
    "AP13" = 12.5
    "AQ13" = 12
    "AR13" = 22
    "H13" = 3.85
    "I13" = 1.62
    "K13" = 2.25
    "L13" = 2.15
    "O13" = 1.21
    "P13" = 3.5
    "W13" = 2.52
    "X13" = 1.4
    "AC17" = 12
    "AF17" = 41
    "AG17" = 29
    "AK17" = 13
    "AM17" = 8.5
    "AQ17" = 15
    "L17" = 2.6
    "M17" = 1.02
    "N17" = 11
    "S17" = 1.83
    "T17" = 1.98
    "AC21" = 7.5
    "AD21" = 11
    "AE21" = 9
    "AF21" = 21
    "AG21" = 19
    "AI21" = 10
    "AM21" = 10
    "AN21" = 17
    "AO21" = 12
    "AP21" = 34
    "AQ21" = 26
    "G21" = 2.2
    "H21" = 3.3
    "I21" = 3.1
    "J21" = 2.88
    "L21" = 3.75
    "S21" = 2
    "T21" = 1.8
    "Y21" = 1.4
    "Z21" = 2.75
    "AP22" = 23
    "G22" = 2.7
    "I22" = 2.5
    "J22" = 3.25
    "L22" = 3.1
    "N22" = 12
    "AA23" = 1.64
    "AB23" = 2.11
    "AC23" = 7.8
    "AD23" = 6.6
    "AE23" = 6.6
    "AF23" = 9.4
    "AI23" = 15
    "AJ23" = 7.2
    "AK23" = 12
    "AL23" = 50
    "AM23" = 16
    "AN23" = 26
    "AO23" = 13
    "AP23" = 80
    "AQ23" = 35
    "AR23" = 35
    "G23" = 1.55
    "H23" = 4
    "I23" = 4.7
    "J23" = 2
    "K23" = 2.47
    "L23" = 5.2
    "O23" = 1.18
    "P23" = 4.5
    "W23" = 2.44
    "X23" = 1.52
    "AA24" = 1.62
    "AB24" = 2.15
    "AC24" = 8.2
    "AD24" = 6.8
    "AE24" = 6.6
    "AF24" = 9.4
    "AG24" = 9.2
    "AI24" = 16
    "AJ24" = 7.4
    "AK24" = 12
    "AL24" = 45
    "AM24" = 16
    "AN24" = 27
    "AO24" = 13
    "AQ24" = 35
    "G24" = 1.54
    "H24" = 4
    "I24" = 4.7
    "J24" = 1.99
    "K24" = 2.49
    "L24" = 5
    "P24" = 4.6
    "W24" = 2.37
    "X24" = 1.54
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
